$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.073.19'
$ws.Range("E2").Value = '  -1.21%  '
$ws.Range("D3").Value = '3.431.35'
$ws.Range("E3").Value = '  -1.36%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '407.56'
$ws.Range("E5").Value = '  -2.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.68'
$ws.Range("E6").Value = '  +0.92%  '
$ws.Range("E7").Value = '  -0.74%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.672'
$ws.Range("E9").Value = '  -3.19%  '
$ws.Range("E10").Value = '  -4.79%  '
$ws.Range("E11").Value = '  -2.38%  '
$ws.Range("E12").Value = '  -1.55%  '
$ws.Range("D13").Value = '3.964.67'
$ws.Range("E13").Value = '  -1.32%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.45'
$ws.Range("E14").Value = '  -3.95%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '19.96'
$ws.Range("E15").Value = '  -1.35%  '
$ws.Range("D16").Value = '3.426.42'
$ws.Range("E16").Value = '  -1.19%  '
$ws.Range("D17").Value = '62.118.45'
$ws.Range("E17").Value = '  -0.95%  '
$ws.Range("E18").Value = '  -2.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.99'
$ws.Range("E19").Value = '  +0.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000131'
$ws.Range("E20").Value = '  -4.91%  '
$ws.Range("E21").Value = '  -5.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '84.69'
$ws.Range("E22").Value = '  +2.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '316.40'
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.81'
$ws.Range("E24").Value = '  -3.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.13'
$ws.Range("E25").Value = '  -2.84%  '
$ws.Range("E26").Value = '  +9.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '29.81'
$ws.Range("E27").Value = '  -3.58%  '
$ws.Range("E28").Value = '  +1.50%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.71'
$ws.Range("E29").Value = '  -0.84%  '
$ws.Range("E30").Value = '  +2.34%  '
$ws.Range("E31").Value = '  -2.94%  '
$ws.Range("E32").Value = '  -3.41%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '42.85'
$ws.Range("E33").Value = '  -2.79%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.41'
$ws.Range("E34").Value = '  -3.81%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0485'
$ws.Range("E36").Value = '  -2.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.89'
$ws.Range("E37").Value = '  -1.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  +0.24%  '
$ws.Range("E39").Value = '  -4.30%  '
$ws.Range("E40").Value = '  -2.08%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.00'
$ws.Range("E41").Value = '  -1.13%  '
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '138.43'
$ws.Range("E42").Value = '  +0.65%  '
$ws.Range("E43").Value = '  -1.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.292'
$ws.Range("E44").Value = '  +1.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.99'
$ws.Range("E45").Value = '  -1.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.81'
$ws.Range("E46").Value = '  -3.17%  '
$ws.Range("E47").Value = '  -1.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '21.45'
$ws.Range("E48").Value = '  -5.69%  '
$ws.Range("D49").Value = '2.135.94'
$ws.Range("E49").Value = '  -4.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.30'
$ws.Range("E50").Value = '  -5.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.88'
$ws.Range("E51").Value = '  +1.10%  '
